# The workbook tracks daily price records for "Espinaca" (spinach) at
# Femacal de La Calera. A new daily record is inserted as row 281, pushing
# all the existing rows (old 281-352) down by one (new 282-353).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 281, shifting rows 281:352 down to 282:353.
$ws.Rows("281:281").Insert()

# Populate the newly inserted row 281 with the new record's data.
$ws.Cells.Item(281, 1).Value = 3
$ws.Cells.Item(281, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(281, 3).Value = "Coquimbo"
$ws.Cells.Item(281, 4).Value = 44722
$ws.Cells.Item(281, 5).Value = 5
$ws.Cells.Item(281, 6).Value = 100112012
$ws.Cells.Item(281, 7).Value = "Espinaca"
$ws.Cells.Item(281, 8).Value = "Sin especificar"
$ws.Cells.Item(281, 9).Value = "Primera"
$ws.Cells.Item(281, 10).Value = 250
$ws.Cells.Item(281, 11).Value = 3500
$ws.Cells.Item(281, 12).Value = 4000
$ws.Cells.Item(281, 13).Value = 3740
$ws.Cells.Item(281, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(281, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(281, 16).Value = 1247
$ws.Cells.Item(281, 17).Value = 3
$ws.Cells.Item(281, 18).Value = "Hortaliza"
